# Quarterly indexing esoteric bug-fix operation
#
# Column A (rows 2-63) holds date serials that were mistakenly written as the
# 1st of the month. The fix re-indexes each of those dates to the 15th of the
# *following* month (a +44/+45 day shift depending on the month length).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $oldSerial = $cell.Value2

    if ($null -eq $oldSerial) { continue }

    $oldDate = [DateTime]::FromOADate($oldSerial)

    $y = $oldDate.Year
    $m = $oldDate.Month + 1
    if ($m -gt 12) {
        $m = $m - 12
        $y = $y + 1
    }

    $newDate = Get-Date -Year $y -Month $m -Day 15 -Hour 0 -Minute 0 -Second 0 -Millisecond 0

    $cell.Value = $newDate.ToOADate()
}
